$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L3").Value = 1.3
$ws.Range("Q3").Value = 1.64
$ws.Range("X3").Value = 25
$ws.Range("F4").Value = 2.12
$ws.Range("G4").Value = 2.78
$ws.Range("H4").Value = 2.88
$ws.Range("J4").Value = 3.15
$ws.Range("K4").Value = 5.8
$ws.Range("N4").Value = 1.66
$ws.Range("P4").Value = 1.66
$ws.Range("Q4").Value = 1.89
$ws.Range("AD7").Value = 9.6
$ws.Range("AI7").Value = 28
$ws.Range("AJ7").Value = 150
$ws.Range("Q8").Value = 1.92
$ws.Range("F10").Value = 1.3
$ws.Range("G10").Value = 1.53
$ws.Range("H10").Value = 1.09
$ws.Range("I10").Value = 46
$ws.Range("P10").Value = 1.25
$ws.Range("Q10").Value = 1.56
$ws.Range("F11").Value = 2.56
$ws.Range("G11").Value = 3.95
$ws.Range("H11").Value = 2.52
$ws.Range("I11").Value = 3.85
$ws.Range("J11").Value = 2.48
$ws.Range("K11").Value = 4.8
$ws.Range("P11").Value = 1.35
$ws.Range("Q11").Value = 2.36
$ws.Range("T12").Value = 2.2
$ws.Range("U12").Value = 1.79
$ws.Range("X12").Value = 11
$ws.Range("U13").Value = 1.93
$ws.Range("X13").Value = 12
$ws.Range("F14").Value = 2.6
$ws.Range("I14").Value = 3.3
$ws.Range("K14").Value = 3.2
$ws.Range("H15").Value = 3.55
$ws.Range("Q15").Value = 1.98
$ws.Range("H16").Value = 8.6
$ws.Range("N16").Value = 3.9
$ws.Range("Q16").Value = 1.88
$ws.Range("Q17").Value = 2.48
$ws.Range("G18").Value = 2.1
$ws.Range("I18").Value = 5
$ws.Range("P18").Value = 1.66
$ws.Range("F19").Value = 2.64
$ws.Range("G19").Value = 2.86
$ws.Range("J19").Value = 3.15
$ws.Range("K19").Value = 3.5
$ws.Range("N19").Value = 3.5
$ws.Range("Q19").Value = 2.02
$ws.Range("F20").Value = 1.88
$ws.Range("G20").Value = 2.42
$ws.Range("H20").Value = 2.9
$ws.Range("I20").Value = 6.2
$ws.Range("J20").Value = 3.3
$ws.Range("K20").Value = 8.199999999999999
$ws.Range("P20").Value = 2.08
$ws.Range("Q20").Value = 1.55
$ws.Range("P21").Value = 3.4
$ws.Range("Q21").Value = 1.26
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 9.4
$ws.Range("H23").Value = 1.42
$ws.Range("I23").Value = 1.46
$ws.Range("K23").Value = 5.7
$ws.Range("N23").Value = 5.3
$ws.Range("P23").Value = 2.48
$ws.Range("R23").Value = 1.58
$ws.Range("T23").Value = 1.84
$ws.Range("U23").Value = 2
$ws.Range("X23").Value = 30
$ws.Range("AC23").Value = 13
$ws.Range("AE23").Value = 15
$ws.Range("AH23").Value = 24
$ws.Range("AJ23").Value = 290
$ws.Range("AL23").Value = 110
$ws.Range("AM23").Value = 130
$ws.Range("AO23").Value = 5.7
$ws.Range("AD25").Value = 34
$ws.Range("AD27").Value = 12
$ws.Range("AK27").Value = 32
$ws.Range("P28").Value = 2.76
$ws.Range("S28").Value = 2.4
$ws.Range("AF28").Value = 8
$ws.Range("O30").Value = 1.44
$ws.Range("U30").Value = 1.91
$ws.Range("F31").Value = 1.45
$ws.Range("G31").Value = 1.83
$ws.Range("H31").Value = 2.2
$ws.Range("J31").Value = 3.25
$ws.Range("P31").Value = 1.48
$ws.Range("Q31").Value = 2.06
$ws.Range("F32").Value = 1.86
$ws.Range("H32").Value = 3.2
$ws.Range("I32").Value = 6
$ws.Range("J32").Value = 3.45
$ws.Range("P32").Value = 1.92
$ws.Range("Q32").Value = 1.67
$ws.Range("H33").Value = 7.4
$ws.Range("N34").Value = 3.05
$ws.Range("F35").Value = 1.92
$ws.Range("S35").Value = 2.9
$ws.Range("T35").Value = 1.69
$ws.Range("AL35").Value = 29
$ws.Range("F36").Value = 1.37
$ws.Range("P36").Value = 2.16
$ws.Range("Q36").Value = 1.52
$ws.Range("F38").Value = 1.09
